$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name "
$ws.Range("B1").Value = "roll no"

$ws.Range("A2").Value = "devyani"
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = "botre"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "botre1"
$ws.Range("B4").Value = 3

$ws.Range("B4").Select()
